$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing DEC_07xx / DEC_08xx test-case identifiers in column A
# (rows 2-19) so they pick up where the previous batch left off: the 18
# scripts that used to be DEC_0795..DEC_0812 become DEC_0815..DEC_0832.
For ($i = 2; $i -le 19; $i++) {
    $newNum = 815 + ($i - 2)
    $ws.Cells.Item($i, 1).Value = "DEC_{0:D4}" -f $newNum
}

# The three rows that used to hold DEC_0813/DEC_0814/DEC_0815 (rows 20-22)
# are emptied out completely -- those scripts are now represented by the
# renumbered rows above.
$ws.Range("A20:J22").ClearContents()

# Five now-unused placeholder rows are removed from the bottom scratch
# area, tightening the sheet from 31 rows down to 26.
$ws.Rows("24:28").Delete()

$ws.Range("C15").Select()
